$d = $word.ActiveDocument

# Locate the "E-post" and "Nettsted" bullet paragraphs that follow the
# "For spørsmål..." contact intro and remove both of them entirely
# (including their paragraph marks), leaving the "GitHub" bullet intact.
$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "E-post:*") {
        $startPara = $p
    }
    if ($t -like "Nettsted:*") {
        $endPara = $p
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $range = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $range.Delete()
}
